$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()
$ws.Range("B2:E24").Value = $true
$ws.Range("A4").Select()
$ws.Range("I7").Select()
